$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '35.114.78'
$ws.Range('E2').Value = '  -0.83%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.897.80'
$ws.Range('E3').Value = '  -0.42%  '
$ws.Range('E4').Value = '  -0.36%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '252.02'
$ws.Range('E5').Value = '  +2.02%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.696'
$ws.Range('E6').Value = '  -0.52%  '
$ws.Range('E7').Value = '  -0.37%  '
$ws.Range('E8').Value = '  +1.88%  '
$ws.Range('E9').Value = '  +0.68%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0754'
$ws.Range('E10').Value = '  +4.02%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0977'
$ws.Range('E11').Value = '  -1.16%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '13.02'
$ws.Range('E12').Value = '  +2.99%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '2.174.70'
$ws.Range('E13').Value = '  -0.31%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.727'
$ws.Range('E14').Value = '  +1.87%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.97'
$ws.Range('E15').Value = '  +1.41%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.884.60'
$ws.Range('E16').Value = '  -1.25%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '35.128.77'
$ws.Range('E17').Value = '  -0.76%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '74.40'
$ws.Range('E18').Value = '  +1.93%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.0₃0837'
$ws.Range('E19').Value = '  +1.37%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '252.01'
$ws.Range('E20').Value = '  +4.30%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.95'
$ws.Range('E21').Value = '  +0.06%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.03'
$ws.Range('E22').Value = '  -1.58%  '
$ws.Range('E23').Value = '  -0.38%  '
$ws.Range('E24').Value = '  +4.70%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.27'
$ws.Range('E25').Value = '  -2.00%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '168.66'
$ws.Range('E26').Value = '  -0.09%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.57'
$ws.Range('E27').Value = '  -1.22%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.42'
$ws.Range('E28').Value = '  -2.61%  '
$ws.Range('E29').Value = '  -1.66%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.128.67'
$ws.Range('E30').Value = '  -0.33%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.32'
$ws.Range('E31').Value = '  +1.90%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.99'
$ws.Range('E32').Value = '  +6.51%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0594'
$ws.Range('E33').Value = '  +3.52%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.62'
$ws.Range('E34').Value = '  +9.69%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.24'
$ws.Range('E35').Value = '  +1.45%  '
$ws.Range('E36').Value = '  -0.36%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.847'
$ws.Range('E37').Value = '  -7.89%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.02'
$ws.Range('E38').Value = '  -0.70%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '17.52'
$ws.Range('E39').Value = '  +5.31%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '98.85'
$ws.Range('E40').Value = '  +0.47%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0213'
$ws.Range('E41').Value = '  +2.36%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0663'
$ws.Range('E42').Value = '  +1.59%  '
$ws.Range('E43').Value = '  -0.66%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.40'
$ws.Range('E44').Value = '  -0.63%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.302.56'
$ws.Range('E45').Value = '  -4.27%  '
$ws.Range('E46').Value = '  +0.07%  '
$ws.Range('E47').Value = '  -1.74%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '6.58'
$ws.Range('E48').Value = '  +0.78%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0768'
$ws.Range('E49').Value = '  +8.25%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '12.09'
$ws.Range('E50').Value = '  -2.24%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '42.78'
$ws.Range('E51').Value = '  -7.25%  '
